# Commit: "fixed to generate log file as csv"
# The underlying edit updates two Item Code values on the
# "Inventory for check stock" sheet: the leading "3" prefix on two item
# codes was corrected to "5" (C4 and C5), and the active selection moved
# to C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory for check stock")

$ws.Range("C4").Value = "5000038608"
$ws.Range("C5").Value = "5000040851"

$ws.Range("C6").Select()
